# "Delete conclusions slide for now"
# Locate the slide titled "Conclusions" (the last slide, #42 - id 589)
# and remove it from the presentation.

$p = $ppt.ActivePresentation

$target = $null
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Conclusions") {
                $target = $slide
                break
            }
        }
    }
    if ($target -ne $null) {
        break
    }
}

if ($target -eq $null) {
    # Fallback: if for some reason the title can't be matched, remove the
    # final slide in the deck, which is where the Conclusions slide lives.
    $target = $p.Slides.Item($p.Slides.Count)
}

$target.Delete()
